$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.573599100112915
$ws.Range("B1").Value = 7.281776428222656
$ws.Range("C1").Value = 6.96187686920166
$ws.Range("D1").Value = 6.307700634002686
$ws.Range("E1").Value = 3.340533256530762
